$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Latest Handback DateTime (col L) and Error Detail (col R) on row 6
$wsZhCn.Range("L6").Value = "2017-02-28 08:11:34"
$wsZhCn.Range("R6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/ea01137e91296f6828c4ef43cd688c9bed35a9af/e2e/10f3806c-2998-43cd-a6a2-45851d8a87dc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/17633b9c133fa9351e548f6c8f6ad97315579176/e2e/10f3806c-2998-43cd-a6a2-45851d8a87dc.md."

# de-de sheet: Latest Handback DateTime (col L) on row 6
$wsDeDe.Range("L6").Value = "2017-02-28 08:11:55"
